$wb = $excel.ActiveWorkbook

# ============================================================
# 1) Insert a new "2022-Q3" worksheet, positioned right after
#    "总计" and before the existing "2022-Q2" sheet.
# ============================================================
$zj  = $wb.Worksheets.Item("总计")
$q2  = $wb.Worksheets.Item("2022-Q2")
$q1  = $wb.Worksheets.Item("2022-Q1")
$q3  = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

# Match page margins used by the other per-quarter sheets
$q3.PageSetup.LeftMargin   = 54
$q3.PageSetup.RightMargin  = 54
$q3.PageSetup.TopMargin    = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Reuse the existing bold+bordered "header" style (B1 on 总计) so the new
# sheet visually matches the other quarter sheets, instead of creating a
# brand-new style entry.
$headerStyleSrc = $zj.Range("B1")
$headerStyleSrc.Copy($q3.Range("B1:H1"))
$headerStyleSrc.Copy($q3.Range("A2:A9"))

# Columns B:G hold text even when the content looks numeric (fund codes,
# percentages, etc. keep their original formatting/leading zeros), so force
# the number format to Text before writing the values.
$q3.Range("B2:G9").NumberFormat = "@"

# --- row 1 ---
$q3.Range("B1").Value = '基金代码'
$q3.Range("C1").Value = '基金名称'
$q3.Range("D1").Value = '基金规模'
$q3.Range("E1").Value = '股票总仓位'
$q3.Range("F1").Value = '仓位占比'
$q3.Range("G1").Value = '持有市值(亿元)'
$q3.Range("H1").Value = '仓位排名'

# --- row 2 ---
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = '014062'
$q3.Range("C2").Value = '景顺长城专精特新量化优选股票A'
$q3.Range("D2").Value = '8.02'
$q3.Range("E2").Value = '91.10'
$q3.Range("F2").Value = '2.48'
$q3.Range("G2").Value = '0.1989'
$q3.Range("H2").Value = 2

# --- row 3 ---
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = '014063'
$q3.Range("C3").Value = '景顺长城专精特新量化优选股票C'
$q3.Range("D3").Value = '5.41'
$q3.Range("E3").Value = '91.10'
$q3.Range("F3").Value = '2.48'
$q3.Range("G3").Value = '0.1342'
$q3.Range("H3").Value = 2

# --- row 4 ---
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = '229002'
$q3.Range("C4").Value = '泰达宏利逆向策略混合'
$q3.Range("D4").Value = '1.59'
$q3.Range("E4").Value = '91.90'
$q3.Range("F4").Value = '1.61'
$q3.Range("G4").Value = '0.0256'
$q3.Range("H4").Value = 9

# --- row 5 ---
$q3.Range("A5").Value = 3
$q3.Range("B5").Value = '001017'
$q3.Range("C5").Value = '泰达宏利改革动力量化策略灵活配置混合A'
$q3.Range("D5").Value = '1.14'
$q3.Range("E5").Value = '91.83'
$q3.Range("F5").Value = '2.11'
$q3.Range("G5").Value = '0.0241'
$q3.Range("H5").Value = 6

# --- row 6 ---
$q3.Range("A6").Value = 4
$q3.Range("B6").Value = '009719'
$q3.Range("C6").Value = '招商增浩一年定期开放混合C'
$q3.Range("D6").Value = '1.34'
$q3.Range("E6").Value = '23.60'
$q3.Range("F6").Value = '0.95'
$q3.Range("G6").Value = '0.0127'
$q3.Range("H6").Value = 5

# --- row 7 ---
$q3.Range("A7").Value = 5
$q3.Range("B7").Value = '161727'
$q3.Range("C7").Value = '招商增荣灵活配置混合（LOF）'
$q3.Range("D7").Value = '0.51'
$q3.Range("E7").Value = '55.32'
$q3.Range("F7").Value = '1.68'
$q3.Range("G7").Value = '0.0086'
$q3.Range("H7").Value = 8

# --- row 8 ---
$q3.Range("A8").Value = 6
$q3.Range("B8").Value = '009718'
$q3.Range("C8").Value = '招商增浩一年定期开放混合A'
$q3.Range("D8").Value = '0.70'
$q3.Range("E8").Value = '23.60'
$q3.Range("F8").Value = '0.95'
$q3.Range("G8").Value = '0.0066'
$q3.Range("H8").Value = 5

# --- row 9 ---
$q3.Range("A9").Value = 7
$q3.Range("B9").Value = '003550'
$q3.Range("C9").Value = '泰达宏利改革动力量化策略灵活配置混合C'
$q3.Range("D9").Value = '0.01'
$q3.Range("E9").Value = '91.83'
$q3.Range("F9").Value = '2.11'
$q3.Range("G9").Value = '0.0002'
$q3.Range("H9").Value = 6

# ============================================================
# 2) Update the "总计" summary sheet: shift the existing two rows
#    down and insert the new 2022-Q3 totals at the top.
# ============================================================
# Give the new A4 cell the same style as the existing index column (A2/A3)
# before writing into it.
$zj.Range("A3").Copy($zj.Range("A4"))

$zj.Range("A4").Value = 2
$zj.Range("B4").Value = '2022-Q1'
$zj.Range("C4").Value = 3
$zj.Range("D4").Value = 0.85

$zj.Range("A3").Value = 1
$zj.Range("B3").Value = '2022-Q2'
$zj.Range("C3").Value = 7
$zj.Range("D3").Value = 0.49

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = '2022-Q3'
$zj.Range("C2").Value = 8
$zj.Range("D2").Value = 0.41

# ============================================================
# 3) Restore "2022-Q1" as the selected/active tab (it was the
#    active sheet before the edit; adding a sheet moves focus).
#    Re-fetch by name: inserting a sheet shifts what index-based
#    worksheet references resolve to, so the earlier $q1 handle
#    can no longer be trusted here.
# ============================================================
$wb.Worksheets.Item("2022-Q1").Activate()
